$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns to row 1
$ws.Range("N1").Value = "Job / Employer"
$ws.Range("O1").Value = "Retirement Date"

# Match the bold header style used by the existing headers (A1:M1)
$ws.Range("N1:O1").Font.Bold = $true

# Set column widths to match target stored widths (col 14 -> 16.1640625, col 15 -> 20.5)
$ws.Columns.Item(14).ColumnWidth = 15.33
$ws.Columns.Item(15).ColumnWidth = 19.67

# Update view state: scroll so column H is the top-left and select P9
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("P9").Select()
